# SOR Testing_Finance.xlsx update
# - Add two new site sheets (Downers Grove Aerospace Illinois, Fort Wayne Indiana)
# - Refresh CVD (turnover) figures across several existing site sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Manila Philippines (sheet style 15)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Manila Philippines")

$ws.Range("G4").Value = 0.0078
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.0076
$ws.Range("J4").Value = 0.0153
$ws.Range("K4").Value = 0.0233

$ws.Range("E5").Value = 0.333333333333333
$ws.Range("E6").Value = 0.333333333333333

$ws.Range("E7").Value = 0.333333333333333
$ws.Range("H7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("L7").Value = 0.3333
$ws.Range("M7").Value = 0.333333333333333
$ws.Range("N7").Value = 0.333333333333333
$ws.Range("O7").Value = 0.333333333333333
$ws.Range("P7").Value = 0.333333333333333
$ws.Range("Q7").Value = 0.333333333333333
$ws.Range("R7").Value = 0.333333333333333
$ws.Range("S7").Value = 0.333333333333333
$ws.Range("T7").Value = 0.333333333333333
$ws.Range("U7").Value = 0.333333333333333
$ws.Range("V7").Value = 0.333333333333333
$ws.Range("W7").Value = 0.333333333333333

$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776

# ---------------------------------------------------------------------------
# 2. Milwaukee Pmc Hq Wisconsin (sheet style 16)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")

$ws.Range("K4").Value = 0.0571

$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0

# ---------------------------------------------------------------------------
# 3. Milwaukee Wisconsin (sheet style 18)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Milwaukee Wisconsin")

# Row 4: switch this row from "PY Actual / Commit-Forecast / Professional Voluntary Turnover"
# placeholder into the Manufacturing Voluntary Turnover / PY Actual row, clearing the
# now-stale monthly cells that no longer carry data.
$ws.Range("D4").Value = "Manufacturing Voluntary Turnover"
$ws.Range("E4").Value = 0.0776
$ws.Range("F4").Value = "PY Actual"
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("P4").ClearContents()
$ws.Range("Q4").ClearContents()
$ws.Range("S4").ClearContents()
$ws.Range("T4").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("V4").ClearContents()

# Row 5: PY Actual -> AOP data source, refreshed ytd figure.
$ws.Range("E5").Value = 0.0776
$ws.Range("F5").Value = "AOP"

# Row 6 (old AOP row) is no longer needed now that row 5 carries the AOP data source.
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 4. Ratingen Germany (sheet style 22)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ratingen Germany")
$ws.Range("L4").Value = 0

# ---------------------------------------------------------------------------
# 5. Rosemont Illinois (sheet style 23)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Rosemont Illinois")
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# New row 10: Commit/Forecast data source for Manufacturing Voluntary Turnover.
$ws.Range("A10").Value = "Corporate"
$ws.Range("B10").Value = "Finance"
$ws.Range("C10").Value = "Rosemont Illinois"
$ws.Range("D10").Value = "Manufacturing Voluntary Turnover"
$ws.Range("E10").Value = 0.0776
$ws.Range("F10").Value = "Commit/Forecast"
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0
$ws.Range("G10:W10").NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# 6. Tipp City Ohio (sheet style 29)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tipp City Ohio")
$ws.Range("L5").ClearContents()

# ---------------------------------------------------------------------------
# 7. Braintree Massachusetts (sheet style 4)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Braintree Massachusetts")
$ws.Range("L4").Value = 0.3333

# ---------------------------------------------------------------------------
# 8. New sheets: Downers Grove Aerospace Illinois, Fort Wayne Indiana
# ---------------------------------------------------------------------------
$headers = @("segment_function","division_function","location","cvd","ytd","data_source", `
  "Jan","Feb","Mar","Q1","Apr","May","Jun","Q2","Jul","Aug","Sep","Q3","Oct","Nov","Dec","Q4","FY")

function Add-SiteSheet([string]$siteName, [string]$tabName) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $tabName

    for ($col = 1; $col -le $headers.Length; $col++) {
        $newSheet.Cells.Item(1, $col).Value = $headers[$col - 1]
    }

    $newSheet.Range("A2").Value = "Corporate"
    $newSheet.Range("B2").Value = "Finance"
    $newSheet.Range("C2").Value = $siteName
    $newSheet.Range("D2").Value = "Manufacturing Voluntary Turnover"
    $newSheet.Range("E2").Value = 0.0776
    $newSheet.Range("F2").Value = "Commit/Forecast"
    $newSheet.Range("G2").ClearContents()
    $newSheet.Range("H2").ClearContents()
    $newSheet.Range("I2").ClearContents()
    $newSheet.Range("J2").ClearContents()
    $newSheet.Range("K2").ClearContents()
    $newSheet.Range("L2").Value = 0
    $newSheet.Range("M2").Value = 0
    $newSheet.Range("N2").Value = 0
    $newSheet.Range("O2").Value = 0
    $newSheet.Range("P2").Value = 0
    $newSheet.Range("Q2").Value = 0
    $newSheet.Range("R2").Value = 0
    $newSheet.Range("S2").Value = 0
    $newSheet.Range("T2").Value = 0
    $newSheet.Range("U2").Value = 0
    $newSheet.Range("V2").Value = 0
    $newSheet.Range("W2").Value = 0
    $newSheet.Range("E2").NumberFormat = "0.0%"
    $newSheet.Range("G2:W2").NumberFormat = "0.0%"

    $newSheet.Range("A1:W1").Columns.AutoFit() | Out-Null
}

Add-SiteSheet "Downers Grove Aerospace Illinois" "Downers Grove Aerospace Illino"
Add-SiteSheet "Fort Wayne Indiana" "Fort Wayne Indiana"
